# "case with 380 kV done" -- refresh the computed per-unit bus-voltage
# results (vm_pu) for the new 380 kV base case.
#
# Columns: A=bus index (unchanged), B:F and I:N hold the per-unit voltages
# that were recomputed for this case; column G is unchanged (flat 1.0) and
# column H has no data (matches the source sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed values for B,C,D,E,F,I,J,K,L,M,N, one row per array element,
# for worksheet rows 2-25 (bus indices 0-23).
$newResults = @(
    (1.02, 1.045620262725187, 1.053264459838836, 1.053472928151599, 1.064602427152493, 1.043599458744019, 1.050679023919927, 1.056010908014736, 1.056218800659506, 1.067317920747931, 1.052171108648462),
    (1.02, 1.046490415189282, 1.053961420501941, 1.054246653795925, 1.06545615015653, 1.043804952033931, 1.051197304170316, 1.056521176234962, 1.056805678612851, 1.067986784283054, 1.052690124916254),
    (1.02, 1.047054150406387, 1.054413010152079, 1.054748297517037, 1.066009683375776, 1.043937073139291, 1.05153268532617, 1.05685128502478, 1.0571857563158, 1.068420052739563, 1.053025982351818),
    (1.02, 1.047291307915095, 1.054603002595822, 1.054959424095585, 1.066242654322417, 1.043992413422732, 1.051673682885591, 1.056990044942881, 1.05734561818232, 1.068602309482655, 1.053167180143946),
    (1.02, 1.047331137230067, 1.054634911565265, 1.054994886937833, 1.066281786687651, 1.04400169335019, 1.051697357154124, 1.05701334226026, 1.057372464178773, 1.068632917635813, 1.053190888032656),
    (1.02, 1.047057318677237, 1.054415548276191, 1.054751117677399, 1.066012795303816, 1.043937813399082, 1.051534569330999, 1.056853139213272, 1.057187892096327, 1.068422487630179, 1.05302786903215),
    (1.02, 1.045914191339385, 1.053499873699783, 1.053734205974415, 1.064890714333233, 1.04366908107328, 1.050854174422633, 1.056183369008945, 1.05641706995229, 1.067543868233833, 1.052346507884971),
    (1.02, 1.04390519703979, 1.051891083946898, 1.051949951188641, 1.062922106553256, 1.043189086200585, 1.049655434719696, 1.055002679714188, 1.055061361228898, 1.065999291607096, 1.051146065834159),
    (1.02, 1.042569553235391, 1.050821853312627, 1.050765715116816, 1.061615626208735, 1.042864792336864, 1.048856487158568, 1.054215317784241, 1.05415937399283, 1.064972132217387, 1.050345983675854),
    (1.02, 1.041992098663444, 1.050359668569289, 1.050254199900313, 1.061051335227377, 1.04272335913023, 1.048510600339803, 1.053874342196069, 1.053769252843058, 1.064527987840425, 1.049999605658133),
    (1.02, 1.041777741316307, 1.050188114440099, 1.050064392546503, 1.060841948272445, 1.042670673225389, 1.048382133281879, 1.053747683612575, 1.053624412821586, 1.064363107982113, 1.049870956162251),
    (1.02, 1.041823715560853, 1.050224907839769, 1.050105098142324, 1.060886852686403, 1.042681981366857, 1.048409689396808, 1.053774852516448, 1.053655478381944, 1.064398470940046, 1.049898551410024),
    (1.02, 1.041974377055317, 1.050345485355098, 1.050238506434285, 1.061034022821966, 1.042719007179613, 1.04849998098689, 1.053863872653661, 1.053757278919145, 1.064514356869684, 1.04998897122452),
    (1.02, 1.042067222509573, 1.050419793343621, 1.050320729206017, 1.061124727859212, 1.042741799975201, 1.048555614072652, 1.053918720266865, 1.053820010689422, 1.064585770671167, 1.050044683315646),
    (1.02, 1.042607895837516, 1.050852543982935, 1.050799689545229, 1.061653106503943, 1.042874157534766, 1.048879443958835, 1.054237946426569, 1.054185274553578, 1.065001621869467, 1.05036897307741),
    (1.02, 1.042947284790055, 1.051124212256328, 1.051100469179636, 1.061984926995172, 1.042956911534865, 1.049082591593516, 1.054438178085122, 1.054414515159027, 1.065262642127073, 1.050572409205033),
    (1.02, 1.043145330197964, 1.051282748752544, 1.05127603095721, 1.062178609502111, 1.043005082869622, 1.049201090140147, 1.054554965616062, 1.05454827019138, 1.065414950772241, 1.050691076033192),
    (1.02, 1.043212873030966, 1.051336818604759, 1.051335913641447, 1.062244673448569, 1.043021491461534, 1.049241496055674, 1.054594786391469, 1.054593884427643, 1.065466894202442, 1.050731539329755),
    (1.02, 1.042910862704689, 1.051095056845387, 1.051068185728127, 1.061949311588494, 1.0429480429116, 1.049060795148318, 1.054416695548676, 1.054389915367346, 1.065234630911565, 1.050550581806384),
    (1.02, 1.041930007261195, 1.050309974903708, 1.05019921569767, 1.060990678911192, 1.042708108171987, 1.048473392058213, 1.053837658571131, 1.053727299300016, 1.064480228712415, 1.049962344536515),
    (1.02, 1.041314086321032, 1.0498170686788, 1.049653972517123, 1.060389198254738, 1.042556376567154, 1.048104131004439, 1.053473566384689, 1.053311082411047, 1.064006457357722, 1.04959255908969),
    (1.02, 1.041640523030367, 1.050078300079188, 1.049942910149151, 1.060707935415359, 1.042636895077625, 1.048299876915802, 1.053666580809011, 1.053531688827273, 1.064257559721081, 1.049788582982699),
    (1.02, 1.04292732001307, 1.051108230682958, 1.051082772853189, 1.06196540423396, 1.042952050562586, 1.0490706440054, 1.054426402595358, 1.054401030822053, 1.065247287789107, 1.050560444649973),
    (1.02, 1.04442392748446, 1.052306421429044, 1.052410303348588, 1.063430002926616, 1.043313936900309, 1.04996530580588, 1.055307963560288, 1.05541153012019, 1.066398157759197, 1.051456376972828)
)

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

for ($i = 0; $i -lt $newResults.Count; $i++) {
    $row = 2 + $i
    $values = $newResults[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $values[$j]
    }
}
